# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback (target content now in sync with en-US).

$wb = $excel.ActiveWorkbook

$statusText   = "Handed back: in sync with en-US"

$md1 = "76379734-f9ef-43c5-bedd-93b772d73204.md"
$md2 = "d636a910-91f2-4d31-8664-869f531d8151.md"

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/47487b203b7f59f99fc7c2eea9a32e9e8f62269a/e2e/76379734-f9ef-43c5-bedd-93b772d73204.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/47487b203b7f59f99fc7c2eea9a32e9e8f62269a/e2e/d636a910-91f2-4d31-8664-869f531d8151.md"

$zhXlf1 = "76379734-f9ef-43c5-bedd-93b772d73204.c9485cc9f3f45d78040a602c7b545f07d3d2b1d5.zh-cn.xlf"
$zhXlf2 = "d636a910-91f2-4d31-8664-869f531d8151.57dcaf765cd2379eb5034e5152e3e4075c8fa5ea.zh-cn.xlf"
$deXlf1 = "76379734-f9ef-43c5-bedd-93b772d73204.c9485cc9f3f45d78040a602c7b545f07d3d2b1d5.de-de.xlf"
$deXlf2 = "d636a910-91f2-4d31-8664-869f531d8151.57dcaf765cd2379eb5034e5152e3e4075c8fa5ea.de-de.xlf"

$zhHandbackDate = "2016-09-04 21:08:46"
$deHandbackDate = "2016-09-04 21:08:53"

# wide-column width (matches existing "wide" columns such as File Name / Path)
$wide = 39.166666666666664
# slightly-narrower width used by the Status columns after the longer text
$statusWidth = 29.166666666666664

# ---------------------------------------------------------------------------
# Overview sheet: refresh the per-language status columns
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = $statusWidth
$overview.Columns.Item(6).ColumnWidth = $statusWidth

# ---------------------------------------------------------------------------
# zh-cn sheet: status + newly-populated "Latest Target File" / handback info
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Hyperlinks.Add($zh.Range("I2"), $url1, "", "", $md1)
$zh.Range("J2").Value = $zhXlf1
$zh.Range("K2").Value = $zhHandbackDate

$zh.Hyperlinks.Add($zh.Range("I3"), $url2, "", "", $md2)
$zh.Range("J3").Value = $zhXlf2
$zh.Range("K3").Value = $zhHandbackDate

$zh.Columns.Item(3).ColumnWidth = $statusWidth
$zh.Columns.Item(9).ColumnWidth = $wide
$zh.Columns.Item(10).ColumnWidth = $wide

# ---------------------------------------------------------------------------
# de-de sheet: status + newly-populated "Latest Target File" / handback info
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Hyperlinks.Add($de.Range("I2"), $url1, "", "", $md1)
$de.Range("J2").Value = $deXlf1
$de.Range("K2").Value = $deHandbackDate

$de.Hyperlinks.Add($de.Range("I3"), $url2, "", "", $md2)
$de.Range("J3").Value = $deXlf2
$de.Range("K3").Value = $deHandbackDate

$de.Columns.Item(3).ColumnWidth = $statusWidth
$de.Columns.Item(9).ColumnWidth = $wide
$de.Columns.Item(10).ColumnWidth = $wide

Write-Output "Report regenerated for handback."
